$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header cell B1: "Typenprüfung" -> "Typenprüfung & Konvertierung" (wrapped, taller row) ---
$ws.Range("B1").Value = "Typenprüfung &`nKonvertierung"
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30.75

# --- Rename the binary/blob conversion function labels in column B ---
# BINARY/VARBINARY/LONGVARBINARY/BLOB/BIT rows used a single generic "bin2hex" label;
# these are now split into more specific conversion function names.
$ws.Range("B18").Value = "base64->hex"   # VARBINARY
$ws.Range("B19").Value = "base64->hex"   # LONGVARBINARY
$ws.Range("B26").Value = "base64->hex"   # BLOB
$ws.Range("B17").Value = "bin->hex"      # BINARY
$ws.Range("B2").Value  = "bit->hex"      # BIT

# --- Restore the cell selection to match the author's saved view ---
$ws.Range("J17").Select()
